$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos table
# with the latest scraped values (GitHub Actions cron update).
# A leading "'" forces cells that look numeric (e.g. "401.45") to stay text,
# matching the original sheet where every Price/Volume cell is a text value;
# the Style reset afterwards clears the quote-prefix formatting so the cell
# keeps its original (default) style.

$ws.Range("D2").Value = '54.515.21'
$ws.Range("E2").Value = '  +5.32%  '
$ws.Range("D3").Value = '3.171.75'
$ws.Range("E3").Value = '  +1.91%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''401.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = '''109.86'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.30%  '
$ws.Range("D8").Value = '''1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +4.19%  '
$ws.Range("D10").Value = '''39.15'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.65%  '
$ws.Range("D11").Value = '''0.0898'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.26%  '
$ws.Range("E12").Value = '  +1.64%  '
$ws.Range("D13").Value = '3.675.93'
$ws.Range("E13").Value = '  +2.08%  '
$ws.Range("D14").Value = '''19.09'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.11%  '
$ws.Range("D15").Value = '''8.07'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.03%  '
$ws.Range("E16").Value = '  +6.53%  '
$ws.Range("D17").Value = '3.174.97'
$ws.Range("E17").Value = '  +2.77%  '
$ws.Range("D18").Value = '''10.58'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.70%  '
$ws.Range("D19").Value = '54.429.18'
$ws.Range("E19").Value = '  +5.01%  '
$ws.Range("D20").Value = '''3.31'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.01%  '
$ws.Range("E21").Value = '  +4.67%  '
$ws.Range("D22").Value = '''12.92'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.62%  '
$ws.Range("D23").Value = '''72.04'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.80%  '
$ws.Range("D24").Value = '''275.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.75%  '
$ws.Range("E25").Value = '  +4.24%  '
$ws.Range("D26").Value = '''8.08'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = '''7.66'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.42%  '
$ws.Range("D28").Value = '''27.81'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.47%  '
$ws.Range("E29").Value = '  +0.27%  '
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("E31").Value = '  +2.26%  '
$ws.Range("D32").Value = '''11.13'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.29%  '
$ws.Range("D33").Value = '''0.0507'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +12.67%  '
$ws.Range("D34").Value = '''36.71'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.73%  '
$ws.Range("E35").Value = '  +1.16%  '
$ws.Range("D36").Value = '''51.31'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.03%  '
$ws.Range("D37").Value = '''3.64'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.17%  '
$ws.Range("D38").Value = '''1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("D39").Value = '''2.88'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +10.19%  '
$ws.Range("E40").Value = '  +10.50%  '
$ws.Range("E41").Value = '  +0.62%  '
$ws.Range("E42").Value = '  +1.72%  '
$ws.Range("D43").Value = '''17.28'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.77%  '
$ws.Range("D44").Value = '''131.99'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.05%  '
$ws.Range("E45").Value = '  +1.12%  '
$ws.Range("D46").Value = '''22.04'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.91%  '
$ws.Range("E47").Value = '  -0.31%  '
$ws.Range("E48").Value = '  -0.53%  '
$ws.Range("D49").Value = '2.101.88'
$ws.Range("E49").Value = '  +2.52%  '
$ws.Range("E50").Value = '  +14.42%  '
$ws.Range("D51").Value = '''0.0339'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.46%  '
